#
# Applies the changes described by the commit "update some bug and template":
#   1. Remove the (now-unnecessary) spellcheck <w:proofErr> wrapper around
#      the "judul_kbli" placeholder token.
#   2. Fill in two previously-blank "Isi bila ada" table cells with their
#      proper merge-field placeholders (jaringan_utilitas / persyaratan_pelaksanaan).
#   3. Drop three stale <w:lastRenderedPageBreak/> markers left over from a
#      previous layout (they're cached rendering hints, not content).
#   4. Merge the two "Surat " / "Keterangan Rencana Kota " footer runs into a
#      single run.
#

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Drop the spellStart/spellEnd proofErr pair wrapping "judul_kbli".
#    There is no direct COM property for <w:proofErr/>, so the paragraph is
#    rebuilt verbatim (same pPr/rPr/rsids) via Range.InsertXML, just without
#    the two proofErr markers.
# ---------------------------------------------------------------------------
$rngKbli = $d.Content
$foundKbli = $rngKbli.Find.Execute('${judul_kbli}', $true, $false, $false, $false, $false,
                                    $true, 0, $false, "", 0)
if ($foundKbli) {
    $xmlKbli = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="53BF3E7E" w14:textId="4118E52A" w:rsidR="00AD1E5F" w:rsidRPr="006E06AA" w:rsidRDefault="006E06AA"><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Bookman Old Style" w:eastAsia="Bookman Old Style" w:hAnsi="Bookman Old Style" w:cs="Bookman Old Style"/><w:lang w:val="en-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Bookman Old Style" w:eastAsia="Bookman Old Style" w:hAnsi="Bookman Old Style" w:cs="Bookman Old Style"/></w:rPr><w:t>${</w:t></w:r><w:r w:rsidRPr="006E06AA"><w:rPr><w:rFonts w:ascii="Bookman Old Style" w:eastAsia="Bookman Old Style" w:hAnsi="Bookman Old Style" w:cs="Bookman Old Style"/><w:lang w:val="en-ID"/></w:rPr><w:t>judul_kbli</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Bookman Old Style" w:eastAsia="Bookman Old Style" w:hAnsi="Bookman Old Style" w:cs="Bookman Old Style"/></w:rPr><w:t>}</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $rngKbli.InsertXML($xmlKbli)
}

# ---------------------------------------------------------------------------
# 2) Fill in the two blank "Isi bila ada" cells (table 1) with placeholders.
#    Both rows' 4th cell starts out reading "Isi bila ada" verbatim; target
#    each cell's own Range so only that occurrence is touched
#    (Replace:=wdReplaceOne keeps the edit scoped to the first hit inside
#    the cell's range).
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)

$cellJaringan = $tbl.Rows.Item(33).Cells.Item(4)
$cellJaringan.Range.Find.Execute("Isi bila ada", $true, $false, $false, $false, $false,
                                  $true, 0, $false, "`${jaringan_utilitas}", 1)

$cellPersyaratan = $tbl.Rows.Item(34).Cells.Item(4)
$cellPersyaratan.Range.Find.Execute("Isi bila ada", $true, $false, $false, $false, $false,
                                     $true, 0, $false, "`${persyaratan_pelaksanaan}", 1)

# ---------------------------------------------------------------------------
# 3) Remove the three stray <w:lastRenderedPageBreak/> markers.
#    Same InsertXML rebuild trick as step 1 -- locate the owning paragraph
#    by its (unique) visible text and re-emit it without the marker.
# ---------------------------------------------------------------------------
function Remove-LastRenderedPageBreak($paragraphText, $xmlBody) {
    # NB: use $d.Content.Paragraphs rather than $d.Paragraphs -- once a
    # Tables.Item(...) anchor has been touched (step 2, above) the bare
    # Document.Paragraphs collection comes back stale/corrupted in this
    # host, while re-deriving the paragraph collection from a fresh
    # Range (Content) keeps reporting the right text/offsets.
    $paras = $d.Content.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text.Trim() -eq $paragraphText) {
            $wrapped = @"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
$xmlBody
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
            $p.Range.InsertXML($wrapped)
            break
        }
    }
}

Remove-LastRenderedPageBreak "Mataram," '<w:p w14:paraId="1E70072B" w14:textId="77777777" w:rsidR="00AD1E5F" w:rsidRDefault="00000000"><w:pPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:left="6094"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Bookman Old Style" w:eastAsia="Bookman Old Style" w:hAnsi="Bookman Old Style" w:cs="Bookman Old Style"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Bookman Old Style" w:eastAsia="Bookman Old Style" w:hAnsi="Bookman Old Style" w:cs="Bookman Old Style"/></w:rPr><w:t>Mataram,</w:t></w:r></w:p>'

Remove-LastRenderedPageBreak "LAMPIRAN KOORDINAT" '<w:p w14:paraId="5FB4D966" w14:textId="77777777" w:rsidR="00AD1E5F" w:rsidRDefault="00000000"><w:pPr><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Bookman Old Style" w:eastAsia="Bookman Old Style" w:hAnsi="Bookman Old Style" w:cs="Bookman Old Style"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Bookman Old Style" w:eastAsia="Bookman Old Style" w:hAnsi="Bookman Old Style" w:cs="Bookman Old Style"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>LAMPIRAN KOORDINAT</w:t></w:r></w:p>'

Remove-LastRenderedPageBreak "LAMPIRAN PETA" '<w:p w14:paraId="46516493" w14:textId="77777777" w:rsidR="00AD1E5F" w:rsidRDefault="00000000"><w:pPr><w:spacing w:after="0"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Bookman Old Style" w:eastAsia="Bookman Old Style" w:hAnsi="Bookman Old Style" w:cs="Bookman Old Style"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Bookman Old Style" w:eastAsia="Bookman Old Style" w:hAnsi="Bookman Old Style" w:cs="Bookman Old Style"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>LAMPIRAN PETA</w:t></w:r></w:p>'

# ---------------------------------------------------------------------------
# 4) Merge the footer's "Surat " + "Keterangan Rencana Kota " runs into one
#    run reading "Surat Keterangan Rencana Kota ". A scoped Find/Replace
#    across the run boundary collapses the matched text into a single run
#    carrying the first run's formatting.
# ---------------------------------------------------------------------------
$sections = $d.Sections
for ($s = 1; $s -le $sections.Count; $s++) {
    $footer = $sections.Item($s).Footers.Item(1)
    if ($footer.Exists) {
        $footer.Range.Find.Execute("Surat Keterangan Rencana Kota ", $true, $false, $false,
                                    $false, $false, $true, 0, $false,
                                    "Surat Keterangan Rencana Kota ", 1)
    }
}
